$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before row 274, shifting the existing rows
# (274..346) down to (276..348). This matches the diff, where the data
# previously at row 274 ends up at row 276, row 275 -> 277, ... row 346 -> 348.
$ws.Rows.Item(274).Resize(2).Insert()

# Populate the two newly-inserted rows with the new records.
# Row 274: Artic Star / Primera
$ws.Range("A274").Value = 7
$ws.Range("B274").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C274").Value = "Ñuble"
$ws.Range("D274").Value = 44900
$ws.Range("E274").Value = 16
$ws.Range("F274").Value = "Fruta"
$ws.Range("G274").Value = 100103
$ws.Range("H274").Value = "Frutos de hueso (carozo)"
$ws.Range("I274").Value = 100103006
$ws.Range("J274").Value = "Nectarín"
$ws.Range("K274").Value = "Artic Star"
$ws.Range("L274").Value = "Primera"
$ws.Range("M274").Value = 120
$ws.Range("N274").Value = 15000
$ws.Range("O274").Value = 16000
$ws.Range("P274").Value = 15500
$ws.Range("Q274").Value = "$/caja 16 kilos empedrada"
$ws.Range("R274").Value = "Región de O'Higgins"
$ws.Range("S274").Value = 969
$ws.Range("T274").Value = 16

# Row 275: Artic Star / Segunda
$ws.Range("A275").Value = 7
$ws.Range("B275").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C275").Value = "Ñuble"
$ws.Range("D275").Value = 44900
$ws.Range("E275").Value = 16
$ws.Range("F275").Value = "Fruta"
$ws.Range("G275").Value = 100103
$ws.Range("H275").Value = "Frutos de hueso (carozo)"
$ws.Range("I275").Value = 100103006
$ws.Range("J275").Value = "Nectarín"
$ws.Range("K275").Value = "Artic Star"
$ws.Range("L275").Value = "Segunda"
$ws.Range("M275").Value = 60
$ws.Range("N275").Value = 14000
$ws.Range("O275").Value = 14000
$ws.Range("P275").Value = 14000
$ws.Range("Q275").Value = "$/caja 16 kilos empedrada"
$ws.Range("R275").Value = "Región de O'Higgins"
$ws.Range("S275").Value = 875
$ws.Range("T275").Value = 16

# Give column D (the date column) the same date number format as the rest
# of the date column so the two new rows render as dates, not raw serials.
$ws.Range("D274:D275").NumberFormat = $ws.Range("D273").NumberFormat
